$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 541.8043
$ws.Range("J17").Value = 541.8043
$ws.Range("L17").Value = 1625.4129
$ws.Range("N17").Value = -1961.4129
$ws.Range("H51").Value = 2099.8572
$ws.Range("J51").Value = 2399.75
$ws.Range("L51").Value = 2399.75
$ws.Range("N51").Value = -3367.75
$ws.Range("H112").Value = 2384.0286
$ws.Range("I112").Value = 683.2222
$ws.Range("J112").Value = 2972.7693
$ws.Range("K112").Value = 2049.6666
$ws.Range("L112").Value = 8918.3079
$ws.Range("M112").Value = -941.6666
$ws.Range("N112").Value = -11134.3079
$ws.Range("H123").Value = 29999
$ws.Range("J123").Value = 29999
$ws.Range("L123").Value = 29999
$ws.Range("N123").Value = -39799

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1645.5385
$ws.Range("I63").Value = 1530.2
$ws.Range("K63").Value = 1530.2
$ws.Range("M63").Value = -844.2
$ws.Range("H66").Value = 1645.5385
$ws.Range("I66").Value = 1530.2
$ws.Range("K66").Value = 7651
$ws.Range("M66").Value = -4219
$ws.Range("H88").Value = 2927
$ws.Range("I88").Value = 2301.6667
$ws.Range("J88").Value = 3083.3333
$ws.Range("K88").Value = 2301.6667
$ws.Range("L88").Value = 3083.3333
$ws.Range("M88").Value = -1895.6667
$ws.Range("N88").Value = -3895.3333
$ws.Range("H91").Value = 2927
$ws.Range("I91").Value = 2301.6667
$ws.Range("J91").Value = 3083.3333
$ws.Range("K91").Value = 2301.6667
$ws.Range("L91").Value = 3083.3333
$ws.Range("M91").Value = -897.6667000000002
$ws.Range("N91").Value = -5891.3333
$ws.Range("H98").Value = 30980
$ws.Range("J98").Value = 30980
$ws.Range("L98").Value = 30980
$ws.Range("N98").Value = -36970
$ws.Range("H122").Value = 1021.8
$ws.Range("I122").Value = 1045.9231
$ws.Range("K122").Value = 3137.7693
$ws.Range("M122").Value = -687.7692999999999
$ws.Range("H132").Value = 1959.5
$ws.Range("I132").Value = 1617.2778
$ws.Range("K132").Value = 4851.8334
$ws.Range("M132").Value = -2321.8334
$ws.Range("H135").Value = 21178.75
$ws.Range("J135").Value = 21178.75
$ws.Range("L135").Value = 21178.75
$ws.Range("N135").Value = -31318.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4583.8
$ws.Range("I134").Value = 1165.0834
$ws.Range("K134").Value = 3495.2502
$ws.Range("M134").Value = -960.2501999999999
$ws.Range("H135").Value = 34994.156
$ws.Range("J135").Value = 34994.156
$ws.Range("L135").Value = 34994.156
$ws.Range("N135").Value = -45134.156

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1217.6608
$ws.Range("I31").Value = 1188.28
$ws.Range("J31").Value = 1462.5
$ws.Range("K31").Value = 1188.28
$ws.Range("L31").Value = 1462.5
$ws.Range("M31").Value = -893.28
$ws.Range("N31").Value = -2052.5
$ws.Range("H34").Value = 1217.6608
$ws.Range("I34").Value = 1188.28
$ws.Range("J34").Value = 1462.5
$ws.Range("K34").Value = 1188.28
$ws.Range("L34").Value = 1462.5
$ws.Range("M34").Value = -986.28
$ws.Range("N34").Value = -1866.5
$ws.Range("H97").Value = 24000
$ws.Range("J97").Value = 24000
$ws.Range("L97").Value = 24000
$ws.Range("N97").Value = -25982
$ws.Range("H105").Value = 772.2
$ws.Range("I105").Value = 744
$ws.Range("J105").Value = 849.75
$ws.Range("K105").Value = 744
$ws.Range("L105").Value = 849.75
$ws.Range("M105").Value = 1003
$ws.Range("N105").Value = -4343.75
$ws.Range("H132").Value = 1790.9354
$ws.Range("I132").Value = 1327.7307
$ws.Range("K132").Value = 3983.1921
$ws.Range("M132").Value = -1453.1921
$ws.Range("H134").Value = 951.46155
$ws.Range("I134").Value = 755.1667
$ws.Range("K134").Value = 2265.5001
$ws.Range("M134").Value = 269.4998999999998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 2000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 2000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 6000
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -6338
$ws.Range("H47").Value = 787.6667
$ws.Range("I47").Value = 481.5
$ws.Range("K47").Value = 1444.5
$ws.Range("M47").Value = -1013.5
$ws.Range("H48").Value = 640
$ws.Range("J48").Value = 980
$ws.Range("L48").Value = 2940
$ws.Range("N48").Value = -3440
$ws.Range("H131").Value = 20003042
$ws.Range("I131").Value = 200000400
$ws.Range("J131").Value = 3336.2
$ws.Range("K131").Value = 600001200
$ws.Range("L131").Value = 10008.6
$ws.Range("M131").Value = -599996160
$ws.Range("N131").Value = -20088.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 253.16667
$ws.Range("I2").Value = 339.66666
$ws.Range("J2").Value = 166.66667
$ws.Range("K2").Value = 339.66666
$ws.Range("L2").Value = 166.66667
$ws.Range("M2").Value = -226.66666
$ws.Range("N2").Value = -392.66667
$ws.Range("H70").Value = 56253264
$ws.Range("I70").Value = 41670216
$ws.Range("K70").Value = 41670216
$ws.Range("M70").Value = -41669946
$ws.Range("H73").Value = 56253264
$ws.Range("I73").Value = 41670216
$ws.Range("K73").Value = 41670216
$ws.Range("M73").Value = -41669280
$ws.Range("H102").Value = 1332.0312
$ws.Range("I102").Value = 1356
$ws.Range("J102").Value = 589
$ws.Range("K102").Value = 1356
$ws.Range("L102").Value = 589
$ws.Range("M102").Value = 266
$ws.Range("N102").Value = -3833
$ws.Range("H113").Value = 1237.1666
$ws.Range("I113").Value = 1157
$ws.Range("J113").Value = 1317.3334
$ws.Range("K113").Value = 1157
$ws.Range("L113").Value = 1317.3334
$ws.Range("M113").Value = 1013
$ws.Range("N113").Value = -5657.3334
$ws.Range("H132").Value = 2015.875
$ws.Range("I132").Value = 1438.3334
$ws.Range("J132").Value = 2882.1875
$ws.Range("K132").Value = 4315.0002
$ws.Range("L132").Value = 8646.5625
$ws.Range("M132").Value = -1785.0002
$ws.Range("N132").Value = -13706.5625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2136.8125
$ws.Range("I82").Value = 2099.2
$ws.Range("J82").Value = 2199.5
$ws.Range("K82").Value = 2099.2
$ws.Range("L82").Value = 2199.5
$ws.Range("M82").Value = -1738.2
$ws.Range("N82").Value = -2921.5
$ws.Range("H85").Value = 2136.8125
$ws.Range("I85").Value = 2099.2
$ws.Range("J85").Value = 2199.5
$ws.Range("K85").Value = 2099.2
$ws.Range("L85").Value = 2199.5
$ws.Range("M85").Value = -851.1999999999998
$ws.Range("N85").Value = -4695.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 14572
$ws.Range("J97").Value = 14572
$ws.Range("L97").Value = 14572
$ws.Range("N97").Value = -16554
$ws.Range("H100").Value = 600
$ws.Range("I100").Value = 600
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1200
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -659
$ws.Range("N100").ClearContents()
$ws.Range("H107").Value = 351.35715
$ws.Range("I107").Value = 356
$ws.Range("K107").Value = 1068
$ws.Range("M107").Value = 852
